$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns G and H
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Area formulas (column G), rows 2-11
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
for ($r = 4; $r -le 11; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 7).Formula = "=(D$r-D$prev)*B$r/100"
}

# Total area (column H), row 2
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Update selection to mirror the post-edit workbook view
$ws.Range("H3").Select()
